# Auto-generated edit script: refresh cryptos list values (rows 2-51)
# per commit 'Updated cryptos list on Wed Mar 27 11:17:31 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the cell to hold $Text verbatim as a string, even when it
    # looks like a number (e.g. '0.620', '12.60', '0.0000306') -- a plain
    # Range.Value assignment would let Excel auto-convert those to numbers
    # and silently drop the significant trailing zero / switch to sci-notation.
    $Range.NumberFormat = '@'
    $Range.Value = $Text
    $Range.ClearFormats()
}

$ws.Range("D2").Value = "70.208.38"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "3.584.76"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "579.43"
$ws.Range("E5").Value = "  -2.29%  "
Set-TextValue $ws.Range("D6") "186.85"
$ws.Range("E6").Value = "  -4.01%  "
$ws.Range("D7").Value = "3.580.93"
$ws.Range("E7").Value = "  -1.30%  "
Set-TextValue $ws.Range("D8") "0.620"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("E9").Value = "  +0.00%  "
Set-TextValue $ws.Range("D10") "0.184"
$ws.Range("E10").Value = "  -1.01%  "
Set-TextValue $ws.Range("D11") "0.652"
$ws.Range("E11").Value = "  -4.02%  "
Set-TextValue $ws.Range("D12") "55.19"
$ws.Range("E12").Value = "  -4.96%  "
Set-TextValue $ws.Range("D13") "0.0000306"
$ws.Range("E13").Value = "  +0.16%  "
Set-TextValue $ws.Range("D14") "9.56"
$ws.Range("E14").Value = "  -4.28%  "
$ws.Range("D15").Value = "4.153.04"
$ws.Range("E15").Value = "  -1.80%  "
Set-TextValue $ws.Range("D16") "19.72"
$ws.Range("E16").Value = "  -3.72%  "
$ws.Range("D17").Value = "3.579.47"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "70.053.73"
$ws.Range("E18").Value = "  -1.09%  "
Set-TextValue $ws.Range("D19") "12.60"
$ws.Range("E19").Value = "  -1.49%  "
Set-TextValue $ws.Range("D20") "0.121"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("E21").Value = "  -2.85%  "
Set-TextValue $ws.Range("D22") "494.23"
$ws.Range("E22").Value = "  +1.14%  "
Set-TextValue $ws.Range("D23") "19.33"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  -5.52%  "
Set-TextValue $ws.Range("D25") "96.81"
$ws.Range("E25").Value = "  +5.89%  "
Set-TextValue $ws.Range("D26") "4.39"
$ws.Range("E26").Value = "  -2.35%  "
Set-TextValue $ws.Range("D27") "11.53"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("E28").Value = "  -6.55%  "
Set-TextValue $ws.Range("D29") "9.36"
$ws.Range("E29").Value = "  -2.25%  "
Set-TextValue $ws.Range("D30") "7.74"
$ws.Range("E30").Value = "  -2.25%  "
Set-TextValue $ws.Range("D31") "31.74"
$ws.Range("E31").Value = "  -3.17%  "
Set-TextValue $ws.Range("D32") "12.14"
$ws.Range("E32").Value = "  -1.07%  "
Set-TextValue $ws.Range("D33") "65.94"
$ws.Range("E33").Value = "  -0.78%  "
Set-TextValue $ws.Range("D34") "0.116"
$ws.Range("E34").Value = "  -6.04%  "
Set-TextValue $ws.Range("D35") "573.82"
$ws.Range("E35").Value = "  -6.96%  "
$ws.Range("E36").Value = "  +15.20%  "
Set-TextValue $ws.Range("D37") "0.414"
$ws.Range("E37").Value = "  +0.57%  "
Set-TextValue $ws.Range("D38") "38.91"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "0.0₃0796"
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D41") "3.49"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D42") "3.17"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  -9.30%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D44") "3.62"
$ws.Range("E44").Value = "  +8.71%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D45") "3.07"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D46") "0.0443"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "3.193.58"
$ws.Range("E47").Value = "  -3.97%  "
Set-TextValue $ws.Range("D48") "9.49"
$ws.Range("E48").Value = "  -1.25%  "
Set-TextValue $ws.Range("D49") "1.57"
$ws.Range("E49").Value = "  +30.65%  "
$ws.Range("E50").Value = "  -2.08%  "
Set-TextValue $ws.Range("D51") "0.999"
$ws.Range("E51").Value = "  -0.21%  "
